$wb = $excel.ActiveWorkbook

# Update "想去人数" (F4, F9) on both "展览" and "全部类型" sheets
$sheetNames = @("展览", "全部类型")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1499
    $ws.Range("F9").Value = 296
}
